# 'Add files via upload' -- refresh of the Daily User Impact Status sheet:
# appends the newly-published daily rows (24-Dec-2025 .. 5-Jan-2026) that
# extend the Date/#Users/#Logged-in/#Errors-bucket table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date(serial), #Total Users, #Logged-in Users, 0/1/2/3-5/6-10/>10 Errors
$rows = @(
    @(46015, 5615, 34, 34, 0, 0, 0, 0, 0),
    @(46016, 5615, 16, 16, 0, 0, 0, 0, 0),
    @(46017, 5615, 3085, 2829, 176, 45, 31, 3, 1),
    @(46018, 5615, 25, 25, 0, 0, 0, 0, 0),
    @(46019, 5615, 29, 28, 1, 0, 0, 0, 0),
    @(46020, 5615, 3495, 3223, 203, 39, 25, 4, 1),
    @(46021, 5615, 3526, 3265, 197, 40, 20, 4, 0),
    @(46022, 5615, 50, 50, 0, 0, 0, 0, 0),
    @(46023, 5615, 39, 38, 1, 0, 0, 0, 0),
    @(46024, 5615, 3159, 2895, 191, 38, 32, 2, 1),
    @(46025, 5615, 49, 49, 0, 0, 0, 0, 0),
    @(46026, 5615, 28, 28, 0, 0, 0, 0, 0),
    @(46027, 5615, 4161, 3862, 216, 48, 31, 4, 0)
)

$data = New-Object "object[,]" $rows.Count,9
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt 9; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}
$ws.Range("A56:I68").Value2 = $data

# Column A keeps the same date number format used by the rows directly above
$ws.Range("A56:A68").NumberFormat = "d-mmm-yy"

# Restore the author's final scroll position / active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 49
$win.ScrollColumn = 1
$ws.Range("J62").Select()
